$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 451.15384
$ws.Range("I12").Value = 388.75
$ws.Range("K12").Value = 388.75
$ws.Range("M12").Value = -218.75

$ws.Range("H41").Value = 1570
$ws.Range("I41").Value = 2359.9
$ws.Range("J41").Value = 441.57144
$ws.Range("K41").Value = 2359.9
$ws.Range("L41").Value = 441.57144
$ws.Range("M41").Value = -1919.9
$ws.Range("N41").Value = -1321.57144

$ws.Range("H92").Value = 33334090
$ws.Range("I92").Value = 41667460
$ws.Range("K92").Value = 41667460
$ws.Range("M92").Value = -41666212

$ws.Range("H96").Value = 52632988
$ws.Range("I96").Value = 919.8461
$ws.Range("J96").Value = 166669140
$ws.Range("K96").Value = 2759.5383
$ws.Range("L96").Value = 500007420
$ws.Range("M96").Value = -1386.5383
$ws.Range("N96").Value = -500010166

$ws.Range("H106").Value = 2803.3489
$ws.Range("I106").Value = 1789.5883
$ws.Range("J106").Value = 6633.1113
$ws.Range("K106").Value = 1789.5883
$ws.Range("L106").Value = 6633.1113
$ws.Range("M106").Value = -1158.5883
$ws.Range("N106").Value = -7895.1113

$ws.Range("H127").Value = 2302686.8
$ws.Range("I127").Value = 1576.091
$ws.Range("J127").Value = 8630741
$ws.Range("K127").Value = 4728.272999999999
$ws.Range("L127").Value = 25892223
$ws.Range("M127").Value = 231.7270000000008
$ws.Range("N127").Value = -25902143

$ws.Range("H129").Value = 506358.53
$ws.Range("I129").Value = 1055.9286
$ws.Range("J129").Value = 1390638.1
$ws.Range("K129").Value = 3167.7858
$ws.Range("L129").Value = 4171914.3
$ws.Range("M129").Value = 1832.2142
$ws.Range("N129").Value = -4181914.3

$ws.Range("H132").Value = 4735.0586
$ws.Range("I132").Value = 3760.1904
$ws.Range("K132").Value = 11280.5712
$ws.Range("M132").Value = -8750.5712

$ws.Range("H138").Value = 11239680
$ws.Range("I138").Value = 58826070
$ws.Range("K138").Value = 176478210
$ws.Range("M138").Value = -176473070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2038.3544
$ws.Range("I32").Value = 1297.6571
$ws.Range("K32").Value = 1297.6571
$ws.Range("M32").Value = -1010.6571

$ws.Range("H61").Value = 1736.5
$ws.Range("I61").Value = 1736.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1736.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1524.5
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 1300.7906
$ws.Range("J74").Value = 1422.5555
$ws.Range("L74").Value = 1422.5555
$ws.Range("N74").Value = -3170.5555

$ws.Range("H77").Value = 1300.7906
$ws.Range("J77").Value = 1422.5555
$ws.Range("L77").Value = 7112.7775
$ws.Range("N77").Value = -15848.7775

$ws.Range("H97").Value = 2999.9614
$ws.Range("I97").Value = 2910.2727
$ws.Range("K97").Value = 2910.2727
$ws.Range("M97").Value = -2414.2727

$ws.Range("H102").Value = 5116.125
$ws.Range("I102").Value = 4042.9092
$ws.Range("K102").Value = 4042.9092
$ws.Range("M102").Value = -2420.9092

$ws.Range("H122").Value = 6794.6504
$ws.Range("I122").Value = 4247.826
$ws.Range("K122").Value = 12743.478
$ws.Range("M122").Value = -10293.478

$ws.Range("H132").Value = 8087.617
$ws.Range("I132").Value = 4286.0527
$ws.Range("K132").Value = 12858.1581
$ws.Range("M132").Value = -10328.1581

$ws.Range("H136").Value = 1736.5
$ws.Range("I136").Value = 1736.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5209.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2659.5
$ws.Range("N136").ClearContents()

$ws.Range("H138").Value = 239744
$ws.Range("J138").Value = 239744
$ws.Range("L138").Value = 239744
$ws.Range("N138").Value = -250024

$ws.Range("H141").Value = 75866.5
$ws.Range("J141").Value = 75866.5
$ws.Range("L141").Value = 75866.5
$ws.Range("N141").Value = -86226.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2381750.8
$ws.Range("I80").Value = 1220.3334
$ws.Range("K80").Value = 1220.3334
$ws.Range("M80").Value = -222.3334

$ws.Range("H83").Value = 2381750.8
$ws.Range("I83").Value = 1220.3334
$ws.Range("K83").Value = 6101.666999999999
$ws.Range("M83").Value = -1109.666999999999

$ws.Range("H94").Value = 3032.7693
$ws.Range("I94").Value = 1383.0952
$ws.Range("K94").Value = 1383.0952
$ws.Range("M94").Value = -932.0952

$ws.Range("H99").Value = 8325.8125
$ws.Range("I99").Value = 8051.6665
$ws.Range("J99").Value = 8678.286
$ws.Range("K99").Value = 8051.6665
$ws.Range("L99").Value = 8678.286
$ws.Range("M99").Value = -6553.6665
$ws.Range("N99").Value = -11674.286

$ws.Range("H105").Value = 2298.389
$ws.Range("I105").Value = 1987.3077
$ws.Range("K105").Value = 1987.3077
$ws.Range("M105").Value = -240.3077000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 484.36
$ws.Range("I7").Value = 169.7
$ws.Range("K7").Value = 169.7
$ws.Range("M7").Value = -56.69999999999999

$ws.Range("H16").Value = 1996.5555
$ws.Range("I16").Value = 1450.0625
$ws.Range("K16").Value = 1450.0625
$ws.Range("M16").Value = -1163.0625

$ws.Range("H22").Value = 1160.8214
$ws.Range("I22").Value = 881.2381
$ws.Range("K22").Value = 881.2381
$ws.Range("M22").Value = -531.2381

$ws.Range("H31").Value = 3091.5557
$ws.Range("I31").Value = 2951.3333
$ws.Range("J31").Value = 3372
$ws.Range("K31").Value = 2951.3333
$ws.Range("L31").Value = 3372
$ws.Range("M31").Value = -2656.3333
$ws.Range("N31").Value = -3962

$ws.Range("H34").Value = 3091.5557
$ws.Range("I34").Value = 2951.3333
$ws.Range("J34").Value = 3372
$ws.Range("K34").Value = 2951.3333
$ws.Range("L34").Value = 3372
$ws.Range("M34").Value = -2749.3333
$ws.Range("N34").Value = -3776

$ws.Range("H92").Value = 44828.715
$ws.Range("J92").Value = 44828.715
$ws.Range("L92").Value = 44828.715
$ws.Range("N92").Value = -49820.715

$ws.Range("H113").Value = 1996.5555
$ws.Range("I113").Value = 1450.0625
$ws.Range("K113").Value = 1450.0625
$ws.Range("M113").Value = 719.9375

$ws.Range("H122").Value = 8930.450000000001
$ws.Range("I122").Value = 9344.857
$ws.Range("K122").Value = 28034.571
$ws.Range("M122").Value = -25584.571

$ws.Range("H132").Value = 4730.2563
$ws.Range("I132").Value = 5143.6875
$ws.Range("K132").Value = 15431.0625
$ws.Range("M132").Value = -12901.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 43.125
$ws.Range("I2").Value = 25.09091
$ws.Range("J2").Value = 65.166664
$ws.Range("K2").Value = 150.54546
$ws.Range("L2").Value = 390.999984
$ws.Range("M2").Value = -37.54545999999999
$ws.Range("N2").Value = -616.999984

$ws.Range("H8").Value = 5399.5
$ws.Range("I8").Value = 5399.5
$ws.Range("K8").Value = 16198.5
$ws.Range("M8").Value = -16059.5

$ws.Range("H45").Value = 2462.7144
$ws.Range("J45").Value = 759.6667
$ws.Range("L45").Value = 2279.0001
$ws.Range("N45").Value = -3343.0001

$ws.Range("H47").Value = 4291
$ws.Range("I47").Value = 10003
$ws.Range("J47").Value = 1435
$ws.Range("K47").Value = 30009
$ws.Range("L47").Value = 4305
$ws.Range("M47").Value = -29578
$ws.Range("N47").Value = -5167

$ws.Range("H107").Value = 1528.0555
$ws.Range("I107").Value = 1200.5555
$ws.Range("J107").Value = 1855.5555
$ws.Range("K107").Value = 3601.6665
$ws.Range("L107").Value = 5566.666499999999
$ws.Range("M107").Value = -1681.6665
$ws.Range("N107").Value = -9406.666499999999

$ws.Range("H119").Value = 6823.2856
$ws.Range("J119").Value = 19932
$ws.Range("L119").Value = 59796
$ws.Range("N119").Value = -69472

$ws.Range("H126").Value = 13757.5
$ws.Range("I126").Value = 14010
$ws.Range("K126").Value = 42030
$ws.Range("M126").Value = -37090

$ws.Range("H134").Value = 10000.75
$ws.Range("I134").Value = 10000.75
$ws.Range("K134").Value = 30002.25
$ws.Range("M134").Value = -24932.25

$ws.Range("H138").Value = 1951.8889
$ws.Range("I138").Value = 2290.923
$ws.Range("J138").Value = 1070.4
$ws.Range("K138").Value = 6872.768999999999
$ws.Range("L138").Value = 3211.2
$ws.Range("M138").Value = -1732.768999999999
$ws.Range("N138").Value = -13491.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18401.428
$ws.Range("I102").Value = 50000
$ws.Range("K102").Value = 50000
$ws.Range("M102").Value = -48378

$ws.Range("H122").Value = 2313.88
$ws.Range("I122").Value = 1767
$ws.Range("J122").Value = 3476
$ws.Range("K122").Value = 5301
$ws.Range("L122").Value = 10428
$ws.Range("M122").Value = -2851
$ws.Range("N122").Value = -15328

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 933
$ws.Range("I100").Value = 562.1539
$ws.Range("K100").Value = 1124.3078
$ws.Range("M100").Value = -583.3078

$ws.Range("H136").Value = 18733.105
$ws.Range("I136").Value = 18536.775
$ws.Range("J136").Value = 19855
$ws.Range("K136").Value = 55610.325
$ws.Range("L136").Value = 59565
$ws.Range("M136").Value = -53060.325
$ws.Range("N136").Value = -64665
